$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the meanEMG / legmaxROM header values (row 1)
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Update CON row (row 2) values
$ws.Range("B2").Value = 452.58925421558689
$ws.Range("C2").Value = 384.50110722503314
$ws.Range("D2").Value = 452.52221118549647
$ws.Range("E2").Value = 381.29154584488936

# Update STR row (row 3) values
$ws.Range("B3").Value = 458.52010153135177
$ws.Range("C3").Value = 387.1899712143674
$ws.Range("D3").Value = 454.25918397107148
$ws.Range("E3").Value = 387.66631878104988

# Update the selection to match the new, narrower selection range
$ws.Range("B1:E3").Select() | Out-Null
